# Apply updated symbol list values to Sheet1 (coinranking crypto snapshot).
# Mirrors the scraped-data refresh: numeric/percentage text cells are re-entered
# as Text (NumberFormat '@') so they stay literal strings like the source feed,
# matching the existing inline-string cell convention used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric/percentage-looking text must be pre-formatted
# as Text, otherwise Excel auto-converts the literal string into a real number.
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11",
    "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16",
    "D17", "E17", "D18", "E18", "D19", "E19", "E20", "E21", "D22", "E22",
    "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "E27", "E28",
    "D40", "E40", "D41", "E41", "E42", "E43", "E44", "D45", "E45", "E46",
    "D47", "E47", "D48", "E48", "E49", "D50", "E50"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Cell value updates (ref = new value)
$ws.Range('D2').Value = '246.80'
$ws.Range('E2').Value = '1.05%'
$ws.Range('D3').Value = '30.21'
$ws.Range('E3').Value = '11.39%'
$ws.Range('D4').Value = '5.174'
$ws.Range('E4').Value = '0.26%'
$ws.Range('D5').Value = '0.05728'
$ws.Range('E5').Value = '1.72%'
$ws.Range('D6').Value = '6.610'
$ws.Range('E6').Value = '2.10%'
$ws.Range('D7').Value = '0.8578'
$ws.Range('E7').Value = '5.07%'
$ws.Range('D8').Value = '0.8804'
$ws.Range('E8').Value = '5.92%'
$ws.Range('B9').Value = 'One'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D9').Value = '0.01032'
$ws.Range('E9').Value = '1,628.91%'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '0.1367'
$ws.Range('E10').Value = '3.01%'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '0.07086'
$ws.Range('E11').Value = '2.70%'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '0.02863'
$ws.Range('E12').Value = '-1.99%'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '0.09388'
$ws.Range('E13').Value = '-0.07%'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '0.001526'
$ws.Range('E14').Value = '1.00%'
$ws.Range('B15').Value = 'CoinExToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D15').Value = '0.04149'
$ws.Range('E15').Value = '-1.92%'
$ws.Range('D16').Value = '0.006046'
$ws.Range('E16').Value = '-1.59%'
$ws.Range('D17').Value = '3.496'
$ws.Range('E17').Value = '-1.75%'
$ws.Range('D18').Value = '3.068'
$ws.Range('E18').Value = '1.69%'
$ws.Range('D19').Value = '2.273'
$ws.Range('E19').Value = '2.09%'
$ws.Range('E20').Value = '2.31%'
$ws.Range('E21').Value = '5.09%'
$ws.Range('D22').Value = '0.1300'
$ws.Range('E22').Value = '0.70%'
$ws.Range('D23').Value = '3.513'
$ws.Range('E23').Value = '-6.47%'
$ws.Range('D24').Value = '0.1380'
$ws.Range('E24').Value = '0.45%'
$ws.Range('D25').Value = '0.001212'
$ws.Range('E25').Value = '-1.08%'
$ws.Range('D26').Value = '0.004499'
$ws.Range('E26').Value = '0.39%'
$ws.Range('E27').Value = '23.48%'
$ws.Range('E28').Value = '-0.01%'
$ws.Range('D40').Value = '0.03785'
$ws.Range('E40').Value = '3.70%'
$ws.Range('D41').Value = '0.005728'
$ws.Range('E41').Value = '66.86%'
$ws.Range('E42').Value = '-22.10%'
$ws.Range('E43').Value = '-14.05%'
$ws.Range('E44').Value = '22.24%'
$ws.Range('D45').Value = '0.00005082'
$ws.Range('E45').Value = '-5.86%'
$ws.Range('E46').Value = '0.02%'
$ws.Range('D47').Value = '0.08897'
$ws.Range('E47').Value = '-18.34%'
$ws.Range('D48').Value = '0.002768'
$ws.Range('E48').Value = '4.75%'
$ws.Range('E49').Value = '0.02%'
$ws.Range('D50').Value = '0.0001999'
$ws.Range('E50').Value = '0.02%'
